# Apply the Dutch (nl-NL) translation workbook update:
#  1. Add a Comment to the "strChkDlgPath" row (row 25) and bump its row
#     height to 30 (the comment now wraps to two lines).
#  2. Insert a brand-new localization row for the "strWindowPos" key right
#     before the "strDlgReset" row (i.e. at worksheet row 32), pushing all
#     the following rows down by one. This also grows the "Tabla13" table
#     by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. strChkDlgPath (row 25): add the missing Comment text -------------
$ws.Range("D25").Value = 'In "settings" form, tab "User interface"'
$ws.Rows.Item(25).RowHeight = 30

# --- 2. Insert the new "strWindowPos" row at row 32 -----------------------
$ws.Rows.Item(32).Insert()

$ws.Range("B32").Value = "localization\strings"
$ws.Range("C32").Value = "strWindowPos"
$ws.Range("D32").Value = 'In "settings" form, tab "User interface"'
$ws.Range("E32").Value = "Remember window position and size on startup"
$ws.Rows.Item(32).RowHeight = 30

# --- 3. Grow the table ("Tabla13") so the new row is included ------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:F204"))
